$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before the old row 13 ("Programa resumido:" row) ---
# This shifts old rows 13-23 down to 15-25, matching the new layout where
# "Docentes responsaveis:" (row 12) now owns two data rows (13 and 14) for
# the two professors, before "Programa resumido:" resumes at row 15.
$ws.Rows("13:14").Insert()

# The insert copies column A's formatting down into the two new rows; the
# target layout has no content in column A for these two rows, so clear it.
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# Borrow the B/C cell formatting (styles 2 and 3) from row 15 (an existing
# data row) so the freshly inserted B13:C14 cells keep the same look as
# every other data row instead of defaulting to column A's style.
$ws.Range("B15:C15").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B15:C15").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)

# --- Fill in the newly inserted "Docentes responsaveis:" rows ---
$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# --- Update "Objetivos:" text (row 10) ---
$objetivos = "Fornecer ao aluno seminários sobre temas atuais de Física, Tecnologia e Engenharia."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Update "Programa resumido:" text (now row 15) ---
$resumido = "Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico."
$ws.Range("B15").Value = $resumido
$ws.Range("C15").Value = $resumido

# --- Update "Programa:" text (now row 17) ---
$programa = "Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços."
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa

# --- Update "Método:" text (now row 20) ---
$metodo = "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina."
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# --- Update "Critério:" text (now row 21) ---
$criterio = "A nota final será calculada pela média aritmética dos relatórios e do seminário."
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# --- Update "Norma de recuperação:" text (now row 22) ---
$norma = "Não há."
$ws.Range("B22").Value = $norma
$ws.Range("C22").Value = $norma

# --- Update "Bibliografia:" text (now row 23) ---
$biblio = "A ser definido de acordo com os temas dos seminários."
$ws.Range("B23").Value = $biblio
$ws.Range("C23").Value = $biblio

Write-Host "edit applied"
